# Post-study competition update: add feature with emotion.
# - Apply a numeric display format to the score columns (B:D)
# - Add new rows documenting additional feature / emotion experiments
# - Update the active selection to the newly added cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Apply the 5-decimal numeric format to the existing score blocks ---
$numFmt = "0.00000_ "
$ws.Range("B1:D5").NumberFormat = $numFmt
$ws.Range("B6:D6").NumberFormat = $numFmt
$ws.Range("B7:D7").NumberFormat = $numFmt
$ws.Range("B9:D10").NumberFormat = $numFmt
$ws.Range("B12:D15").NumberFormat = $numFmt

# --- New rows under the existing data (row 12 gains a note; rows 13-15 add more runs) ---
# (shared-string entries are created in the same order the author typed them)
$ws.Range("F13").Value = "增加1个属性，包含多个得分最高单词"

$ws.Range("B15").Value = 0.91309689999999999
$ws.Range("C15").Value = 0.92390000000000005
$ws.Range("D15").Value = 0.89522999999999997
$ws.Range("F15").Value = "增加多个属性，包含多个得分最高单词"

$ws.Range("F12").Value = "增加了一个属性，关于问号和问词"

$ws.Range("F14").Value = "补充newsDesk&sectionName, 增加1个属性，包含多个得分最高单词"

# --- New block at row 17: emotion.csv feature ---
$ws.Range("A17").Value = "emotion.csv"
$ws.Range("B17").NumberFormat = $numFmt
$ws.Range("B17").Value = 0.92717749999999999
$ws.Range("C17").Value = 0.92286999999999997
$ws.Range("D17").Value = 0.89654999999999996
$ws.Range("F17").Value = "整段文字的emotion"

# --- Update selection to match the last edited cell ---
$ws.Range("E17").Select()
